$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet from Chinese to English
$ws.Name = "Performance Analysis"

# Title row
$ws.Range("A1").Value = "Performance Analysis Report: qwen3 (EXTEND)"

# Header row (row 3)
$ws.Range("A3").Value = "Operator Name"
$ws.Range("B3").Value = "Type"
$ws.Range("H3").Value = "Input"
$ws.Range("I3").Value = "Output"
$ws.Range("J3").Value = "Weight"
$ws.Range("K3").Value = "Compute(us)"
$ws.Range("L3").Value = "Memory(us)"
$ws.Range("M3").Value = "Transfer(us)"
$ws.Range("N3").Value = "Single Layer Latency(us)"
$ws.Range("O3").Value = "Total Time(ms)"
$ws.Range("P3").Value = "Percent(%)"
$ws.Range("Q3").Value = "Weight/Single GPU All Layers"

# Summary labels rows 15-25
$ws.Range("A15").Value = "Compute Time (ms)"
$ws.Range("A16").Value = "Memory Time (ms)"
$ws.Range("A17").Value = "Transfer Time (ms)"
$ws.Range("A18").Value = "Total Time (ms)"
$ws.Range("A21").Value = "Performance Bottleneck"
$ws.Range("B21").Value = "dense_gate_up_proj (Total Time: 68.719 ms)"
$ws.Range("A24").Value = "Throughput TPS"
$ws.Range("A25").Value = "Weight Memory/Single GPU (GB)"

# New row 26: KV Cache Memory/Single GPU (GB)
# Clone formatting from row 25's cells (A25/B25) first so the new row matches
# the existing "label + value" style pairing without minting stray new fonts.
$ws.Range("A25").Copy()
$ws.Range("A26").PasteSpecial(-4122)
$ws.Range("B25").Copy()
$ws.Range("B26").PasteSpecial(-4122)

$ws.Range("A26").Value = "KV Cache Memory/Single GPU (GB)"
$ws.Range("B26").Value = 0.001221
$ws.Range("B26").NumberFormat = "0.000000"

$excel.CutCopyMode = $false
